$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark from the "All ... team member"
#    bullet (it moves to the "Will Comber" edit below).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) "Will comber" -> "Will Comber", split across three runs with a
#    collapsed "_GoBack" bookmark marking the last edit point (right after
#    the capitalised "C"), matching Word's own last-edit-position bookmark.
$range = $d.Content
$range.Find.Execute("Will comber", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $d.Range($range.Start, $range.End)
$target.Text = ""

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Will </w:t></w:r><w:r><w:t>C</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>omber</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint = $d.Range($target.Start, $target.Start)
$insertPoint.InsertXML($xml)
